$d = $word.ActiveDocument
$found = $d.Content.Find.Execute("possible: make a system like in cross code for accuracy, an upgrade that increases the accuracy", $true, $false, $false, $false, $false, $true, 1, $false, "accuracy time", 2)
Write-Output $found
